$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header style from E1 (existing header) to F1:H1
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Boolean outlier flag data for rows 2-25, columns F (KNN), G (SVM), H (RF)
$flags = @(
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,0,0),
    @(1,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,1,0),
    @(0,0,0)
)

for ($i = 0; $i -lt $flags.Length; $i++) {
    $row = $i + 2
    $vals = $flags[$i]
    $ws.Cells.Item($row, 6).Value = [bool]$vals[0]
    $ws.Cells.Item($row, 7).Value = [bool]$vals[1]
    $ws.Cells.Item($row, 8).Value = [bool]$vals[2]
}
